# "fix ar and name excel column"
# Remove the 10 stale/duplicate product rows from the import sheet.
# Deleting from the bottom-most row upward keeps the remaining row
# numbers stable while each EntireRow delete is applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(100, 98, 94, 49, 24, 22, 19, 16, 15, 10)

foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
